$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.695.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.98%  "

$ws.Range("D3").Value = "'1.870.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.31%  "

$ws.Range("E4").Value = "  -0.83%  "

$ws.Range("D5").Value = "'246.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.35%  "

$ws.Range("D6").Value = "'0.686"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.19%  "

$ws.Range("E7").Value = "  -0.87%  "

$ws.Range("D8").Value = "'42.02"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.38%  "

$ws.Range("D9").Value = "'0.345"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.00%  "

$ws.Range("D10").Value = "'51.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.43%  "

$ws.Range("D11").Value = "'0.0731"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.54%  "

$ws.Range("E12").Value = "  -2.79%  "

$ws.Range("D13").Value = "'2.146.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.15%  "

$ws.Range("D14").Value = "'12.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.88%  "

$ws.Range("D15").Value = "'0.710"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.82%  "

$ws.Range("D16").Value = "'4.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.71%  "

$ws.Range("D17").Value = "'1.887.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.41%  "

$ws.Range("D18").Value = "'34.714.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.95%  "

$ws.Range("D19").Value = "'72.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.80%  "

$ws.Range("D20").Value = "'0.0₃0814"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.96%  "

$ws.Range("D21").Value = "'243.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.63%  "

$ws.Range("D22").Value = "'12.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.80%  "

$ws.Range("D23").Value = "'4.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.68%  "

$ws.Range("E24").Value = "  -1.03%  "

$ws.Range("D25").Value = "'2.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.64%  "

$ws.Range("E26").Value = "  -4.85%  "

$ws.Range("D27").Value = "'164.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.86%  "

$ws.Range("D28").Value = "'8.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.82%  "

$ws.Range("D29").Value = "'18.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.29%  "

$ws.Range("E30").Value = "  -5.33%  "

$ws.Range("D31").Value = "'4.128.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.07%  "

$ws.Range("D32").Value = "'1.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.82%  "

$ws.Range("E33").Value = "  -1.03%  "

$ws.Range("D34").Value = "'0.0574"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.11%  "

$ws.Range("D35").Value = "'4.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.62%  "

$ws.Range("E36").Value = "  -0.99%  "

$ws.Range("D37").Value = "'0.823"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.91%  "

$ws.Range("D38").Value = "'1.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -18.77%  "

$ws.Range("D39").Value = "'1.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.64%  "

$ws.Range("D40").Value = "'97.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.53%  "

$ws.Range("D41").Value = "'16.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.38%  "

$ws.Range("D42").Value = "'0.0660"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.04%  "

$ws.Range("D43").Value = "'0.0209"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.29%  "

$ws.Range("D44").Value = "'1.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.25%  "

$ws.Range("D45").Value = "'1.283.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.79%  "

$ws.Range("D46").Value = "'2.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.42%  "

$ws.Range("E47").Value = "  -0.99%  "

$ws.Range("D48").Value = "'0.0781"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.56%  "

$ws.Range("D49").Value = "'2.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.25%  "

$ws.Range("D50").Value = "'12.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.94%  "

$ws.Range("D51").Value = "'6.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.72%  "
